# Fruta / hortaliza, semanal
#
# Two new weekly observations are inserted into the "Mango" sheet right
# above the existing row 77 record, pushing the former rows 77-113 down
# to rows 79-115 (dimension grows from A1:T113 to A1:T115).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 77 (Excel shifts row 77.. down to 79..,
# copying the formatting - including the date number-format on column D -
# from the row directly above, same as native Excel "Insert Copied/Above" behaviour).
$ws.Rows.Item(77).Resize(2).Insert()

# ---- New row 77 -------------------------------------------------------
$ws.Cells.Item(77, 1).Value = 4
$ws.Cells.Item(77, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(77, 3).Value = "Los Lagos"
$ws.Cells.Item(77, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(77, 5).Value = 10
$ws.Cells.Item(77, 6).Value = "Fruta"
$ws.Cells.Item(77, 7).Value = 100108
$ws.Cells.Item(77, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(77, 9).Value = 100108002
$ws.Cells.Item(77, 10).Value = "Mango"
$ws.Cells.Item(77, 11).Value = "Sin especificar"
$ws.Cells.Item(77, 12).Value = "Primera"
$ws.Cells.Item(77, 13).Value = 300
$ws.Cells.Item(77, 14).Value = 7500
$ws.Cells.Item(77, 15).Value = 8000
$ws.Cells.Item(77, 16).Value = 7750
$ws.Cells.Item(77, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(77, 18).Value = "Perú"
$ws.Cells.Item(77, 19).Value = 1938
$ws.Cells.Item(77, 20).Value = 4

# ---- New row 78 -------------------------------------------------------
$ws.Cells.Item(78, 1).Value = 4
$ws.Cells.Item(78, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(78, 3).Value = "Los Lagos"
$ws.Cells.Item(78, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(78, 5).Value = 10
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100108
$ws.Cells.Item(78, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(78, 9).Value = 100108002
$ws.Cells.Item(78, 10).Value = "Mango"
$ws.Cells.Item(78, 11).Value = "Sin especificar"
$ws.Cells.Item(78, 12).Value = "Segunda"
$ws.Cells.Item(78, 13).Value = 100
$ws.Cells.Item(78, 14).Value = 5000
$ws.Cells.Item(78, 15).Value = 5000
$ws.Cells.Item(78, 16).Value = 5000
$ws.Cells.Item(78, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(78, 18).Value = "Perú"
$ws.Cells.Item(78, 19).Value = 1250
$ws.Cells.Item(78, 20).Value = 4
